# Fix for "reading double value from excel" — adds three new test rows to
# the "Number" sheet (normal-precision doubles, including the classic
# 0.1 + 0.2 floating point case and a large integral double), plus a new
# empty "Text" worksheet at the end of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Number")

# Row 8: calculate -> 0.1 + 0.2 (computed), expected text "0.3"
# Expected-value column (C) is entered before the label column (A) so the
# shared-string table fills in the same order as the authored workbook.
$ws.Range("C8").Value = "'0.3"
$ws.Range("A8").Value = "calculate"
$ws.Range("B8").Formula = "= 0.1 + 0.2"

# Row 9: double -> 1.332, expected text "1.332"
$ws.Range("C9").Value = "'1.332"
$ws.Range("A9").Value = "double"
$ws.Range("B9").Value = 1.332

# Row 10: bigvalue -> 200000000000, expected text "200000000000"
# Here the label (A) is entered before the expected text (C).
$ws.Range("A10").Value = "bigvalue"
$ws.Range("C10").Value = "'200000000000"
$ws.Range("B10").Value = 200000000000

# Match the authored page setup for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the cursor to A11 (just past the new data) and keep "Number" active.
$ws.Activate()
[void]$ws.Range("A11").Select()

# Add a new, empty "Text" worksheet as the last tab in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
[void]$wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "Text"

# Re-activate "Number" so it remains the selected/visible tab.
$ws.Activate()
